# Update countries & provincias Spain
#
# 1) Swap ranking of Israel / Polonia: Polonia's updated totals now
#    exceed Israel's, so Polonia moves up to row 33 and Israel drops to row 34.
# 2) Update the "Datos actualizados..." timestamp string.
# 3) Refresh the numeric stats for Rumania (row 39), Indonesia (row 40)
#    and Sri Lanka (row 105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33 / 34: Polonia overtakes Israel ---------------------------------
$ws.Cells.Item(33, 1).Value = "Polonia"
$ws.Cells.Item(33, 2).Value = 16561
$ws.Cells.Item(33, 3).Value = 235
$ws.Cells.Item(33, 4).Value = 6131
$ws.Cells.Item(33, 5).Value = 9603
$ws.Cells.Item(33, 6).Value = 160
$ws.Cells.Item(33, 7).Value = 16
$ws.Cells.Item(33, 8).Value = 827

$ws.Cells.Item(34, 1).Value = "Israel"
$ws.Cells.Item(34, 2).Value = 16506
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 11843
$ws.Cells.Item(34, 5).Value = 4405
$ws.Cells.Item(34, 6).Value = 66
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 258

# --- Row 39: Rumania ---------------------------------------------------------
$ws.Cells.Item(39, 5).Value = 7352
$ws.Cells.Item(39, 7).Value = 9
$ws.Cells.Item(39, 8).Value = 991

# --- Row 40: Indonesia -------------------------------------------------------
$ws.Cells.Item(40, 2).Value = 14749
$ws.Cells.Item(40, 3).Value = 484
$ws.Cells.Item(40, 4).Value = 3063
$ws.Cells.Item(40, 5).Value = 10679
$ws.Cells.Item(40, 7).Value = 16
$ws.Cells.Item(40, 8).Value = 1007

# --- Row 105: Sri Lanka -------------------------------------------------------
$ws.Cells.Item(105, 4).Value = 366
$ws.Cells.Item(105, 5).Value = 494

# --- Title timestamp ----------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 10:35"
